$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.513.01'
$ws.Range('E2').Value = '  +5.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.725.94'
$ws.Range('E3').Value = '  +4.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.43'
$ws.Range('E5').Value = '  +3.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5355'
$ws.Range('E6').Value = '  +2.98%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2663'
$ws.Range('E8').Value = '  +1.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06610'
$ws.Range('E9').Value = '  +4.96%  '
$ws.Range('E10').Value = '  +6.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07688'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.602'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.724.62'
$ws.Range('E13').Value = '  +4.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.964.33'
$ws.Range('E14').Value = '  +4.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5804'
$ws.Range('E15').Value = '  +4.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8299'
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('E17').Value = '  +4.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '27.524.05'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '218.83'
$ws.Range('E19').Value = '  +12.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.723'
$ws.Range('E21').Value = '  +3.01%  '
$ws.Range('E22').Value = '  +1.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.038'
$ws.Range('E23').Value = '  +2.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.59'
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.747'
$ws.Range('E26').Value = '  +13.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1236'
$ws.Range('E27').Value = '  +4.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.342'
$ws.Range('E28').Value = '  +2.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.51'
$ws.Range('E29').Value = '  +4.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05491'
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.301'
$ws.Range('E31').Value = '  +2.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.551'
$ws.Range('E32').Value = '  +3.43%  '
$ws.Range('E33').Value = '  +3.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.664'
$ws.Range('E34').Value = '  +6.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.859'
$ws.Range('E35').Value = '  +2.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9579'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.426'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5939'
$ws.Range('E38').Value = '  +6.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01650'
$ws.Range('E39').Value = '  +5.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.892'
$ws.Range('E40').Value = '  +2.80%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8484'
$ws.Range('E41').Value = '  +3.33%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.046.63'
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.003'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.27'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.870.55'
$ws.Range('E45').Value = '  +4.76%  '
$ws.Range('E46').Value = '  +2.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.83'
$ws.Range('E47').Value = '  +2.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4485'
$ws.Range('E48').Value = '  +3.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.204'
$ws.Range('E49').Value = '  +4.07%  '
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05253'
$ws.Range('E51').Value = '  +3.02%  '
